$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.02258322285507441
$ws.Range("C2").Value = 9.226618575922256
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 71517.89157740913
$ws.Range("G2").Value = 71527.85351205892
